$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-07-03 Thursday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-07-04 Friday", 2) | Out-Null
$d.Content.Find.Execute("134÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "249÷6=", 2) | Out-Null
$d.Content.Find.Execute("738÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "771÷6=", 2) | Out-Null
$d.Content.Find.Execute("143÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "177÷4=", 2) | Out-Null
$d.Content.Find.Execute("870÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "661÷3=", 2) | Out-Null
$d.Content.Find.Execute("173÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "451÷5=", 2) | Out-Null
$d.Content.Find.Execute("734÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "606÷4=", 2) | Out-Null
$d.Content.Find.Execute("117÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "677÷8=", 2) | Out-Null
$d.Content.Find.Execute("692÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "442÷7=", 2) | Out-Null
$d.Content.Find.Execute("377÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "587÷9=", 2) | Out-Null
$d.Content.Find.Execute("457÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "820÷8=", 2) | Out-Null
$d.Content.Find.Execute("146÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "645÷5=", 2) | Out-Null
$d.Content.Find.Execute("578÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "258÷8=", 2) | Out-Null
$d.Content.Find.Execute("845÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "749÷5=", 2) | Out-Null
$d.Content.Find.Execute("575÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "344÷7=", 2) | Out-Null
$d.Content.Find.Execute("230÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "715÷3=", 2) | Out-Null
$d.Content.Find.Execute("430÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "774÷4=", 2) | Out-Null
$d.Content.Find.Execute("657÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "677÷6=", 2) | Out-Null
$d.Content.Find.Execute("961÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "808÷9=", 2) | Out-Null
$d.Content.Find.Execute("490÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "843÷8=", 2) | Out-Null
$d.Content.Find.Execute("631÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "745÷9=", 2) | Out-Null
$d.Content.Find.Execute("659÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "924÷8=", 2) | Out-Null
$d.Content.Find.Execute("130÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "923÷7=", 2) | Out-Null
$d.Content.Find.Execute("742÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "377÷3=", 2) | Out-Null
$d.Content.Find.Execute("920÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "281÷4=", 2) | Out-Null
$d.Content.Find.Execute("520÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "468÷4=", 2) | Out-Null

Write-Host "Replacements complete"